$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 data: title1 / "I have change this description right now" / "bi bi-tiktok"
# becomes: "Mohammad Ali Jarjoumah" / "Full Stack Web Developer" / "bi bi-linkedin"
$ws.Range("A2").Value = "Mohammad Ali Jarjoumah"
$ws.Range("B2").Value = "Full Stack Web Developer"
$ws.Range("C2").Value = "bi bi-linkedin"

# Update the active selection to E3 (was D9)
$ws.Range("E3").Select()
